{"js": "// 1. Update the letter date.\nconst dateHits = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"September 21, 2025\", \"Replace\");\n  await context.sync();\n}\n\n// 2. Split the mailing address line \"7155 Almaden Place, San Jose CA 95120\"\n//    into two separate paragraphs: \"7155 Almaden Place\" and \"San Jose, CA 95120\".\n//    (The same sentence also appears inside the summary table further down \u2014\n//    only the first, body-level occurrence is targeted by this edit.)\nconst addrHits = context.document.body.search(\"7155 Almaden Place, San Jose CA 95120\", { matchCase: true });\naddrHits.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < addrHits.items.length; i++) {\n  addrHits.items[i].parentTableOrNullObject.load(\"isNullObject\");\n}\nawait context.sync();\nconst addrRange = addrHits.items.find((r) => r.parentTableOrNullObject.isNullObject);\nif (addrRange) {\n  const addrParagraph = addrRange.paragraphs.getFirst();\n  // Shrink the existing paragraph down to just the street address.\n  addrRange.insertText(\"7155 Almaden Place\", \"Replace\");\n  // Insert a new paragraph right after it carrying the city/state/zip,\n  // matching the same paragraph/run formatting as the original line.\n  const cityParagraph = addrParagraph.insertParagraph(\"San Jose, CA 95120\", \"After\");\n  cityParagraph.font.set({ name: \"Arial\", size: 11 });\n  await context.sync();\n}\n\n// 3. Remove the now-superfluous empty \"No Spacing\" paragraph that sat\n//    directly below \"...Board of Directors\".\nconst boardHits = context.document.body.search(\"Board of Directors\", { matchCase: true });\nboardHits.load(\"items\");\nawait context.sync();\nif (boardHits.items.length > 0) {\n  const boardParagraph = boardHits.items[0].paragraphs.getFirst();\n  const nextParagraph = boardParagraph.getNext();\n  nextParagraph.load(\"text\");\n  await context.sync();\n  if (nextParagraph.text.trim() === \"\") {\n    nextParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the letter date.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"September 19, 2025\") {\n        $p.Range.Text = \"September 21, 2025\"\n        break\n    }\n}\n\n# 2. Split the mailing address line \"7155 Almaden Place, San Jose CA 95120\"\n#    into two separate paragraphs: \"7155 Almaden Place\" and \"San Jose, CA 95120\".\n#    (The same sentence also appears inside the summary table further down \u2014\n#    only the first, body-level occurrence is targeted by this edit.)\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"7155 Almaden Place, San Jose CA 95120\") {\n        $p.Range.Text = \"7155 Almaden Place`rSan Jose, CA 95120\"\n        break\n    }\n}\n\n# 3. Remove the now-superfluous empty \"No Spacing\" paragraph that sat\n#    directly below \"...Board of Directors\".\n$boardPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Board of Directors*\") {\n        $boardPara = $p\n        break\n    }\n}\nif ($boardPara -ne $null) {\n    $nextPara = $boardPara.Next()\n    if ($nextPara -ne $null -and $nextPara.Range.Text.Trim() -eq \"\" -and $nextPara.Style.NameLocal -eq \"No Spacing\") {\n        $nextPara.Range.Delete()\n    }\n}\n"}
